$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # brp_ingeschrevenpersoon
$ws2 = $wb.Worksheets.Item(2)   # brp_nationaliteit
$ws3 = $wb.Worksheets.Item(3)   # brp_partners
$ws4 = $wb.Worksheets.Item(4)   # brp_ouders
$ws5 = $wb.Worksheets.Item(5)   # brp_kinderen
$ws6 = $wb.Worksheets.Item(6)   # brp_reisdocument

# --- sheet1 (brp_ingeschrevenpersoon): onvolledige_datum fix ---
# AS42 (geboorte__datum) becomes an incomplete date "00000000" entered as text
$ws1.Range("AS42").Value = "'00000000"

# sheet3 header T1
$ws3.Range("T1").Value = "indicatie_onjuist"

# sheet3 existing rows T2..T37 = NULL
$ws3.Range("T2").Value = "NULL"
$ws3.Range("T3").Value = "NULL"
$ws3.Range("T4").Value = "NULL"
$ws3.Range("T5").Value = "NULL"
$ws3.Range("T6").Value = "NULL"
$ws3.Range("T7").Value = "NULL"
$ws3.Range("T8").Value = "NULL"
$ws3.Range("T9").Value = "NULL"
$ws3.Range("T10").Value = "NULL"
$ws3.Range("T11").Value = "NULL"
$ws3.Range("T12").Value = "NULL"
$ws3.Range("T13").Value = "NULL"
$ws3.Range("T14").Value = "NULL"
$ws3.Range("T15").Value = "NULL"
$ws3.Range("T16").Value = "NULL"
$ws3.Range("T17").Value = "NULL"
$ws3.Range("T18").Value = "NULL"
$ws3.Range("T19").Value = "NULL"
$ws3.Range("T20").Value = "NULL"
$ws3.Range("T21").Value = "NULL"
$ws3.Range("T22").Value = "NULL"
$ws3.Range("T23").Value = "NULL"
$ws3.Range("T24").Value = "NULL"
$ws3.Range("T25").Value = "NULL"
$ws3.Range("T26").Value = "NULL"
$ws3.Range("T27").Value = "NULL"
$ws3.Range("T28").Value = "NULL"
$ws3.Range("T29").Value = "NULL"
$ws3.Range("T30").Value = "NULL"
$ws3.Range("T31").Value = "NULL"
$ws3.Range("T32").Value = "NULL"
$ws3.Range("T33").Value = "NULL"
$ws3.Range("T34").Value = "NULL"
$ws3.Range("T35").Value = "NULL"
$ws3.Range("T36").Value = "NULL"
$ws3.Range("T37").Value = "NULL"

# sheet3 row38
$ws3.Range("A38").Value = "33809a77-390a-4acd-9598-6833888495a9"
$ws3.Range("B38").Value = 43
$ws3.Range("C38").Value = "NULL"
$ws3.Range("D38").Value = "H"
$ws3.Range("E38").Value = "M"
$ws3.Range("F38").Value = "Ben"
$ws3.Range("G38").Value = "Niet"
$ws3.Range("H38").Value = "NULL"
$ws3.Range("I38").Value = "NULL"
$ws3.Range("J38").Value = "Luik"
$ws3.Range("K38").Value = 19911111
$ws3.Range("L38").Value = 5010
$ws3.Range("M38").Value = "NULL"
$ws3.Range("N38").Value = "NULL"
$ws3.Range("O38").Value = "NULL"
$ws3.Range("P38").Value = "NULL"
$ws3.Range("Q38").Value = "s-Gravenhage"
$ws3.Range("R38").Value = 20110426
$ws3.Range("S38").Value = 6030
$ws3.Range("T38").Value = "O"

# sheet3 row39
$ws3.Range("A39").Value = "6e0b6bc4-3344-4ec1-b87c-37099e48d4fe"
$ws3.Range("B39").Value = 24
$ws3.Range("C39").Value = "NULL"
$ws3.Range("D39").Value = "H"
$ws3.Range("E39").Value = "M"
$ws3.Range("F39").Value = "Leonel"
$ws3.Range("G39").Value = "Nada"
$ws3.Range("H39").Value = "de"
$ws3.Range("I39").Value = "NULL"
$ws3.Range("J39").Value = "Barcelona"
$ws3.Range("K39").Value = 19790730
$ws3.Range("L39").Value = 6037
$ws3.Range("M39").Value = "Groningen"
$ws3.Range("N39").Value = "N"
$ws3.Range("O39").Value = 20060714
$ws3.Range("P39").Value = 6030
$ws3.Range("Q39").Value = "Groningen"
$ws3.Range("R39").Value = 20050126
$ws3.Range("S39").Value = 6030
$ws3.Range("T39").Value = "NULL"

# sheet4 header N1
$ws4.Range("N1").Value = "indicatie_onjuist"

# sheet4 existing rows N2..N13 = NULL
$ws4.Range("N2").Value = "NULL"
$ws4.Range("N3").Value = "NULL"
$ws4.Range("N4").Value = "NULL"
$ws4.Range("N5").Value = "NULL"
$ws4.Range("N6").Value = "NULL"
$ws4.Range("N7").Value = "NULL"
$ws4.Range("N8").Value = "NULL"
$ws4.Range("N9").Value = "NULL"
$ws4.Range("N10").Value = "NULL"
$ws4.Range("N11").Value = "NULL"
$ws4.Range("N12").Value = "NULL"
$ws4.Range("N13").Value = "NULL"

# sheet4 row14
$ws4.Range("A14").Value = "906fd71b-ab43-4f9d-852a-9898ac8da836"
$ws4.Range("B14").Value = 2
$ws4.Range("D14").Value = "NULL"
$ws4.Range("E14").Value = 1
$ws4.Range("F14").Value = 19830526
$ws4.Range("G14").Value = "Ali"
$ws4.Range("H14").Value = "Baba"
$ws4.Range("I14").Value = "NULL"
$ws4.Range("J14").Value = "NULL"
$ws4.Range("K14").Value = "Riyad"
$ws4.Range("L14").Value = 19560904
$ws4.Range("M14").Value = 5018
$ws4.Range("N14").Value = "NULL"

# sheet4 row15
$ws4.Range("A15").Value = "8a328c17-3305-4a8b-938b-52581f7ebeda"
$ws4.Range("B15").Value = 2
$ws4.Range("C15").Value = 999999370
$ws4.Range("D15").Value = "NULL"
$ws4.Range("E15").Value = 2
$ws4.Range("F15").Value = 19830526
$ws4.Range("G15").Value = "NULL"
$ws4.Range("H15").Value = "NULL"
$ws4.Range("I15").Value = "NULL"
$ws4.Range("J15").Value = "NULL"
$ws4.Range("K15").Value = "NULL"
$ws4.Range("L15").Value = "NULL"
$ws4.Range("M15").Value = "NULL"
$ws4.Range("N15").Value = "O"

# sheet5 header M1
$ws5.Range("M1").Value = "indicatie_onjuist"

# sheet5 existing rows M2..M11 = NULL
$ws5.Range("M2").Value = "NULL"
$ws5.Range("M3").Value = "NULL"
$ws5.Range("M4").Value = "NULL"
$ws5.Range("M5").Value = "NULL"
$ws5.Range("M6").Value = "NULL"
$ws5.Range("M7").Value = "NULL"
$ws5.Range("M8").Value = "NULL"
$ws5.Range("M9").Value = "NULL"
$ws5.Range("M10").Value = "NULL"
$ws5.Range("M11").Value = "NULL"

# sheet5 row12
$ws5.Range("A12").Value = "14f15ad4-5836-419c-a578-a450a61920b6"
$ws5.Range("B12").Value = 39
$ws5.Range("C12").Value = "NULL"
$ws5.Range("D12").Value = "NULL"
$ws5.Range("E12").Value = 20190117
$ws5.Range("F12").Value = "Jesse"
$ws5.Range("G12").Value = "Groenen"
$ws5.Range("H12").Value = "NULL"
$ws5.Range("I12").Value = "NULL"
$ws5.Range("J12").Value = "Istanbul"
$ws5.Range("K12").Value = 20190114
$ws5.Range("L12").Value = 6043
$ws5.Range("M12").Value = "NULL"

# sheet5 row13
$ws5.Range("A13").Value = "c66674be-5789-4ede-9066-b85d32520d93"
$ws5.Range("B13").Value = 46
$ws5.Range("C13").Value = 999999023
$ws5.Range("D13").Value = "NULL"
$ws5.Range("E13").Value = 19830526
$ws5.Range("F13").Value = "NULL"
$ws5.Range("G13").Value = "NULL"
$ws5.Range("H13").Value = "NULL"
$ws5.Range("I13").Value = "NULL"
$ws5.Range("J13").Value = "NULL"
$ws5.Range("K13").Value = "NULL"
$ws5.Range("L13").Value = "NULL"
$ws5.Range("M13").Value = "O"


# --- View / selection changes ---

# sheet3 (brp_partners): new selection
$ws3.Activate()
$ws3.Range("B32").Select()

# sheet4 (brp_ouders): new selection (range B11:B13)
$ws4.Activate()
$ws4.Range("B11:B13").Select()

# sheet5 (brp_kinderen): new selection
$ws5.Activate()
$ws5.Range("C13").Select()

# sheet6 (brp_reisdocument): keep selection, no longer the active/selected tab
$ws6.Activate()
$ws6.Range("H3").Select()

# sheet1 (brp_ingeschrevenpersoon): becomes the active tab, zoom + selection change
$ws1.Activate()
$excel.ActiveWindow.Zoom = 125
$ws1.Range("N3").Select()
